$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.280011666666667
$ws.Range("H2").Value = 9.840035
$ws.Range("I2").Value = 0.5749674887425973
$ws.Range("J2").Value = 0.5749674887425973
$ws.Range("M2").Value = 0.05057900000000001
$ws.Range("N2").Value = 0.151737
$ws.Range("O2").Value = 0.01400296657613869
$ws.Range("P2").Value = 0.01400296657613869
$ws.Range("Q2").Value = 0.1658997100883333
$ws.Range("R2").Value = 1.493097390795
$ws.Range("S2").Value = 0.008051250527228989
$ws.Range("T2").Value = 0.008051250527228989
$ws.Range("G3").Value = 3.280011666666667
$ws.Range("H3").Value = 9.840035
$ws.Range("I3").Value = 0.5749674887425973
$ws.Range("J3").Value = 0.5749674887425973
$ws.Range("O3").Value = 0.146324388539341
$ws.Range("P3").Value = 0.146324388539341
$ws.Range("Q3").Value = 1.733573632811111
$ws.Range("R3").Value = 15.6021626953
$ws.Range("S3").Value = 0.08413176622026097
$ws.Range("T3").Value = 0.08413176622026097
$ws.Range("G4").Value = 3.280011666666667
$ws.Range("H4").Value = 9.840035
$ws.Range("I4").Value = 0.5749674887425973
$ws.Range("J4").Value = 0.5749674887425973
$ws.Range("O4").Value = 0.8396726448845202
$ws.Range("P4").Value = 0.8396726448845202
$ws.Range("Q4").Value = 9.947995490671111
$ws.Range("R4").Value = 89.53195941604
$ws.Range("S4").Value = 0.4827844719951073
$ws.Range("T4").Value = 0.4827844719951073
$ws.Range("I5").Value = 0.1205821735470086
$ws.Range("J5").Value = 0.1205821735470086
$ws.Range("M5").Value = 0.05057900000000001
$ws.Range("N5").Value = 0.151737
$ws.Range("O5").Value = 0.01400296657613869
$ws.Range("P5").Value = 0.01400296657613869
$ws.Range("Q5").Value = 0.034792484836
$ws.Range("R5").Value = 0.313132363524
$ws.Range("S5").Value = 0.001688508145856916
$ws.Range("T5").Value = 0.001688508145856916
$ws.Range("I6").Value = 0.1205821735470086
$ws.Range("J6").Value = 0.1205821735470086
$ws.Range("O6").Value = 0.146324388539341
$ws.Range("P6").Value = 0.146324388539341
$ws.Range("S6").Value = 0.01764411281301072
$ws.Range("T6").Value = 0.01764411281301073
$ws.Range("I7").Value = 0.1205821735470086
$ws.Range("J7").Value = 0.1205821735470086
$ws.Range("O7").Value = 0.8396726448845202
$ws.Range("P7").Value = 0.8396726448845202
$ws.Range("S7").Value = 0.1012495525881409
$ws.Range("T7").Value = 0.1012495525881409
$ws.Range("I8").Value = 0.3044503377103941
$ws.Range("J8").Value = 0.3044503377103941
$ws.Range("M8").Value = 0.05057900000000001
$ws.Range("N8").Value = 0.151737
$ws.Range("O8").Value = 0.01400296657613869
$ws.Range("P8").Value = 0.01400296657613869
$ws.Range("Q8").Value = 0.08784535430500001
$ws.Range("R8").Value = 0.7906081887450002
$ws.Range("S8").Value = 0.004263207903052785
$ws.Range("T8").Value = 0.004263207903052786
$ws.Range("I9").Value = 0.3044503377103941
$ws.Range("J9").Value = 0.3044503377103941
$ws.Range("O9").Value = 0.146324388539341
$ws.Range("P9").Value = 0.146324388539341
$ws.Range("R9").Value = 8.2614822483
$ws.Range("S9").Value = 0.04454850950606928
$ws.Range("T9").Value = 0.04454850950606929
$ws.Range("I10").Value = 0.3044503377103941
$ws.Range("J10").Value = 0.3044503377103941
$ws.Range("O10").Value = 0.8396726448845202
$ws.Range("P10").Value = 0.8396726448845202
$ws.Range("S10").Value = 0.255638620301272
$ws.Range("T10").Value = 0.255638620301272
